$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 5340
$ws.Range("F6").Value = 621
$ws.Range("F7").Value = 464
$ws.Range("F8").Value = 248
$ws.Range("F9").Value = 1119
$ws.Range("F11").Value = 167
$ws.Range("F13").Value = 753
$ws.Range("F14").Value = 391
$ws.Range("F16").Value = 98
$ws.Range("F18").Value = 7
$ws.Range("F19").Value = 380
$ws.Range("F20").Value = 6179
$ws.Range("F21").Value = 50
$ws.Range("F22").Value = 53
$ws.Range("F24").Value = 7159
$ws.Range("F27").Value = 3277
$ws.Range("F28").Value = 396
$ws.Range("F29").Value = 785
$ws.Range("F30").Value = 4471
$ws.Range("F31").Value = 330
$ws.Range("F34").Value = 1223
$ws.Range("F35").Value = 115
$ws.Range("F36").Value = 34
$ws.Range("F37").Value = 9
$ws.Range("F38").Value = 970
$ws.Range("F39").Value = 1241
$ws.Range("F41").Value = 6

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1165

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1165
$ws.Range("F8").Value = 5340
$ws.Range("F9").Value = 621
$ws.Range("F10").Value = 464
$ws.Range("F11").Value = 248
$ws.Range("F12").Value = 1119
$ws.Range("F14").Value = 167
$ws.Range("F16").Value = 753
$ws.Range("F17").Value = 391
$ws.Range("F20").Value = 98
$ws.Range("F22").Value = 7
$ws.Range("F23").Value = 380
$ws.Range("F24").Value = 6179
$ws.Range("F25").Value = 6179
$ws.Range("F26").Value = 50
$ws.Range("F27").Value = 53
$ws.Range("F29").Value = 7159
$ws.Range("F32").Value = 3277
$ws.Range("F33").Value = 396
$ws.Range("F34").Value = 785
$ws.Range("F35").Value = 4471
$ws.Range("F36").Value = 330
$ws.Range("F40").Value = 1223
$ws.Range("F41").Value = 115
$ws.Range("F42").Value = 34
$ws.Range("F43").Value = 9
$ws.Range("F44").Value = 970
$ws.Range("F45").Value = 1241
$ws.Range("F48").Value = 6
